$wb = $excel.ActiveWorkbook

# "Repayment Schedule" is the sheet that needs a new column inserted
# before column N (shifting N:P -> O:Q), and becomes the active/selected sheet.
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N; existing N/O/P data shifts to O/P/Q.
$ws.Columns("N").Insert()

# Update the active cell / selection on this sheet.
$ws.Range("R8").Select()

# Make "Repayment Schedule" the active (selected) sheet/tab.
$ws.Activate()

$wb.Save()
